# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.722.64"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.95"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.99%  "

$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0685"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.046.52"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.11"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +11.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.798.50"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.720.06"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.632"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("E17").Value = "  +3.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.60"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.58"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0780"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.48"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.85%  "

$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.76"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0515"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.438.96"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.75%  "

$ws.Range("E36").Value = "  -1.21%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.634"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.52%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0189"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "82.98"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("E40").Value = "  +4.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.904"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0503"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.94"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.943.36"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.15"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.74"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.30%  "
